$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 157
$startDevice = 3000176

for ($i = 0; $i -lt 5; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = 10002
    $ws.Cells.Item($r, 2).Value = $startDevice + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

$ws.Range("B157").Select()
$excel.ActiveWindow.ScrollRow = 152
$excel.ActiveWindow.ScrollColumn = 1
